# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
#    and populate it with the Q1-2022 fund holdings table.
# 2. Prepend a new "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet positioned after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$newSheet.Name = "2022-Q1"

# Header row (bold, bordered, centered/top-aligned - matches the other
# sheets' header styling)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # starts at column B
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# Data rows: index, code, name, size, position, ratio, value, rank
$rows = @(
    @(0, "213003", "宝盈策略增长混合",       "10.28", "94.38", "6.53", "0.6713", 6),
    @(1, "213002", "宝盈泛沿海增长混合",     "5.10",  "93.76", "6.86", "0.3499", 3),
    @(2, "000965", "汇丰晋信新动力混合",     "0.99",  "91.64", "4.11", "0.0407", 4),
    @(3, "540004", "汇丰晋信2026周期混合",   "1.14",  "31.29", "1.72", "0.0196", 7),
    @(4, "005104", "富荣福康混合A",          "0.08",  "87.88", "3.10", "0.0025", 3),
    @(5, "005105", "富荣福康混合C",          "0.04",  "87.88", "3.10", "0.0012", 3)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $newSheet.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
    $aCell.Borders.Weight = 2

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 7).Style = "Normal"
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet - insert a new 2022-Q1 row on top
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Existing rows (before edit): row2=2021-Q4(6,1.92) row3=2021-Q3(5,2.24)
# New layout: row2=2022-Q1(6,1.09) row3=2021-Q4(6,1.92) row4=2021-Q3(5,2.24)
$totals = @(
    @("2022-Q1", 6, 1.09),
    @("2021-Q4", 6, 1.92),
    @("2021-Q3", 5, 2.24)
)

$r = 2
foreach ($t in $totals) {
    $aCell = $zj.Cells.Item($r, 1)
    $aCell.Value = $r - 2
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
    $aCell.Borders.Weight = 2

    $zj.Cells.Item($r, 2).Value = $t[0]
    $zj.Cells.Item($r, 3).Value = $t[1]
    $zj.Cells.Item($r, 4).Value = $t[2]
    $r = $r + 1
}

# Restore the original active sheet/tab (creating the sheet above makes it
# active by default - put the selection back on the first sheet).
$wb.Worksheets.Item("2021-Q3").Activate()

